$wb = $excel.ActiveWorkbook

$gVal = [double]"5.68631262647114e+23"

# --- Sheet 1: ROW50-FE-LIFTER -> add row 25 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$ws1.Cells.Item(25,1).Value2 = 45736.63332814815
$ws1.Cells.Item(25,2).Value2 = "0x01,0x90"
$ws1.Cells.Item(25,3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item(25,4).Value2 = "0x01,0x7a"
$ws1.Cells.Item(25,5).Value2 = "0xe"
$ws1.Cells.Item(25,6).Value2 = 400
$ws1.Cells.Item(25,7).Value2 = $gVal
$ws1.Cells.Item(25,8).Value2 = 378
$ws1.Cells.Item(25,9).Value2 = 14
$ws1.Cells.Item(25,1).NumberFormat = $ws1.Cells.Item(24,1).NumberFormat

# --- Sheet 2: ROW50-MID-LIFTER -> add row 27 ---
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$ws2.Cells.Item(27,1).Value2 = 45736.61004629629
$ws2.Cells.Item(27,2).Value2 = "0x01,0x90 "
$ws2.Cells.Item(27,3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item(27,4).Value2 = "0x01,0x82"
$ws2.Cells.Item(27,5).Value2 = "0x19"
$ws2.Cells.Item(27,6).Value2 = 400
$ws2.Cells.Item(27,7).NumberFormat = "@"
$ws2.Cells.Item(27,7).Value2 = "568631262647113771663628"
$ws2.Cells.Item(27,7).NumberFormat = "General"
$ws2.Cells.Item(27,8).Value2 = 386
$ws2.Cells.Item(27,9).Value2 = 25
$ws2.Cells.Item(27,1).NumberFormat = $ws2.Cells.Item(26,1).NumberFormat

# --- Sheet 3: ROW11-FE-LIFTER -> add row 25 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$ws3.Cells.Item(25,1).Value2 = 45736.65663960648
$ws3.Cells.Item(25,2).Value2 = "0x01,0x90"
$ws3.Cells.Item(25,3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item(25,4).Value2 = "0x01,0x7a"
$ws3.Cells.Item(25,5).Value2 = "0x14"
$ws3.Cells.Item(25,6).Value2 = 400
$ws3.Cells.Item(25,7).Value2 = $gVal
$ws3.Cells.Item(25,8).Value2 = 378
$ws3.Cells.Item(25,9).Value2 = 20
$ws3.Cells.Item(25,1).NumberFormat = $ws3.Cells.Item(24,1).NumberFormat

# --- Sheet 4: ROW11-MID-LIFTER -> add row 25 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$ws4.Cells.Item(25,1).Value2 = 45736.8038628125
$ws4.Cells.Item(25,2).Value2 = "0x01,0x90"
$ws4.Cells.Item(25,3).Value2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item(25,4).Value2 = "0x01,0x82"
$ws4.Cells.Item(25,5).Value2 = "0x19"
$ws4.Cells.Item(25,6).Value2 = 400
$ws4.Cells.Item(25,7).Value2 = $gVal
$ws4.Cells.Item(25,8).Value2 = 386
$ws4.Cells.Item(25,9).Value2 = 25
$ws4.Cells.Item(25,1).NumberFormat = $ws4.Cells.Item(24,1).NumberFormat

Write-Output "done"
